# Apply content edits to the "Metrics Sheet" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics Sheet")

# B3: clarify report type description (add "un retard de rapport")
$ws.Range("B3").Value = "Définit si les données envoyées sont un rapport, un retard de rapport ou un rapport de rappel "

# B4: pluralize "rapport" -> "rapports"
$ws.Range("B4").Value = "Nombre de rapports "

# D5: extend comparison formula with year comparison
$ws.Range("D5").Value = " Comparaison(date_rapport_de_rappel(le mois) & mois_de_rapport) && Comparaison(date_rapport_de_rappel(l'année) & annee_en_cours)"

# A7: rename metric label
$ws.Range("A7").Value = "Nombre des agents assidus"

# B7: update description of the metric
$ws.Range("B7").Value = "Calcule le nombre de rapports total délivrés par les agents du début de l'année jusqu'au mois en cours"

# D8: update calculation description
$ws.Range("D8").Value = "date_rapport - 1er jour succédant le mois_de_rapport"

# Update the active selection to C3 (matches the edited workbook's saved selection)
$ws.Range("C3").Select()
